$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.883656666666667
$ws.Range("H2").Value = 5.650970000000001
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 6.743825
$ws.Range("N2").Value = 20.231475
$ws.Range("O2").Value = 0.06175298453338176
$ws.Range("P2").Value = 0.06809253940846915
$ws.Range("Q2").Value = 12.70305092008334
$ws.Range("R2").Value = 114.32745828075
$ws.Range("S2").Value = 0.06175298453338176
$ws.Range("T2").Value = 0.06809253940846915

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.883656666666667
$ws.Range("H3").Value = 5.650970000000001
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 69.657865
$ws.Range("N3").Value = 208.973595
$ws.Range("O3").Value = 0.6378547871531949
$ws.Range("P3").Value = 0.7033368922862506
$ws.Range("Q3").Value = 131.2115017930167
$ws.Range("R3").Value = 1180.90351613715
$ws.Range("S3").Value = 0.6378547871531949
$ws.Range("T3").Value = 0.7033368922862506

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.883656666666667
$ws.Range("H4").Value = 5.650970000000001
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.263827666666667
$ws.Range("N4").Value = 3.791483
$ws.Range("O4").Value = 0.01157282852869501
$ws.Range("P4").Value = 0.01276089388411081
$ws.Range("Q4").Value = 2.380617409834445
$ws.Range("R4").Value = 21.42555668851
$ws.Range("S4").Value = 0.01157282852869501
$ws.Range("T4").Value = 0.01276089388411081

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.883656666666667
$ws.Range("H5").Value = 5.650970000000001
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.038909
$ws.Range("N5").Value = 3.116727
$ws.Range("O5").Value = 0.009513255668495419
$ws.Range("P5").Value = 0.01048988549144043
$ws.Range("Q5").Value = 1.95694786391
$ws.Range("R5").Value = 17.61253077519
$ws.Range("S5").Value = 0.009513255668495419
$ws.Range("T5").Value = 0.01048988549144043

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.883656666666667
$ws.Range("H6").Value = 5.650970000000001
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 30.502036
$ws.Range("N6").Value = 61.004072
$ws.Range("O6").Value = 0.2793061441162328
$ws.Range("P6").Value = 0.205319788929729
$ws.Range("Q6").Value = 57.45536345830668
$ws.Range("R6").Value = 344.73218074984
$ws.Range("S6").Value = 0.2793061441162328
$ws.Range("T6").Value = 0.205319788929729
